$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")
Write-Host $ws.Name
